# Update auto-generated "Viện phí" (hospital fee) records in rows 2 and 3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 1844
$ws.Range("E2").Value = 46200602537
$ws.Range("X2").Value = "DN4127460129048"
$ws.Range("AV2").Value = 1094172

# Row 3
$ws.Range("A3").Value = 1845
$ws.Range("E3").Value = 46200602538
$ws.Range("X3").Value = "DN4127460129049"
$ws.Range("AV3").Value = 1094172
